# Apply "contingencies with rene fine" edit to the lines_states worksheet.
#
# Rows 8 and 9 (previously the first two "extr" entries) are renamed to
# line7 / line8, all subsequent "extr" rows shift down two rows, and two
# brand-new "extr7" / "extr8" rows are appended at the bottom. Several
# from_bus / to_bus / in_service values are also refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# name, from_bus, to_bus, in_service for rows 2..17 (row index in col A is name-1)
$data = @(
    @(2,  0, "line1",  7,  9, $true),
    @(3,  1, "line2",  9,  8, $false),
    @(4,  2, "line3",  8, 10, $true),
    @(5,  3, "line4",  8, 11, $true),
    @(6,  4, "line5", 10,  5, $true),
    @(7,  5, "line6", 12,  8, $true),
    @(8,  6, "line7", 14, 11, $true),
    @(9,  7, "line8", 16,  9, $true),
    @(10, 8, "extr1",  5, 12, $true),
    @(11, 9, "extr2",  5,  9, $true),
    @(12, 10, "extr3", 10, 11, $false),
    @(13, 11, "extr4",  7,  8, $false),
    @(14, 12, "extr5",  9, 11, $false),
    @(15, 13, "extr6",  7, 11, $false),
    @(16, 14, "extr7",  5,  7, $true),
    @(17, 15, "extr8",  8,  5, $false)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

# Rows 16 and 17 are brand new - give column A the same formatting
# (bold font, thin box border, centered) used by the other "name" cells.
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16:A17").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0
